# Updated cryptos list on Mon May 22 05:22:38 UTC 2023 with GitHub Actions
#
# Updates the Price (D) / Volume(1h) (E) columns for each coin row, and
# also fixes two rows whose Coin/Link/Price/Volume data had been swapped
# (WrappedBTC <-> Avalanche at rows 20/21, Quant <-> EnergySwap at rows
# 47/48).
#
# Price values in column D are stored as plain text in the workbook (e.g.
# "308.69", "0.07550", "26.905.70") rather than as numbers, so every write
# to column D is entered with a leading apostrophe to force Excel to keep
# it as text instead of auto-converting it to a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'26.933.03"
$ws.Range("E2").Value = "  -1.54%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.817.31"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.30%  "

# Row 5 - USDC
$ws.Range("E5").Value = "  -0.20%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'308.78"
$ws.Range("E6").Value = "  -1.65%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4634"
$ws.Range("E7").Value = "  -2.26%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3647"
$ws.Range("E8").Value = "  -1.22%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.07215"
$ws.Range("E9").Value = "  -3.32%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "'0.8587"
$ws.Range("E10").Value = "  -3.12%  "

# Row 11 - Solana
$ws.Range("D11").Value = "'19.75"
$ws.Range("E11").Value = "  -3.48%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.07561"
$ws.Range("E12").Value = "  +3.11%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "'1.794.29"
$ws.Range("E13").Value = "  -5.66%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.319"
$ws.Range("E14").Value = "  -2.42%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'91.72"
$ws.Range("E15").Value = "  -1.50%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "'6.467"
$ws.Range("E16").Value = "  -1.74%  "

# Row 17 - BinanceUSD
$ws.Range("E17").Value = "  +0.06%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.000008609"
$ws.Range("E18").Value = "  -2.38%  "

# Row 19 - Dai
$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "  -0.36%  "

# Row 20 - was WrappedBTC, now Avalanche
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'14.43"
$ws.Range("E20").Value = "  -2.45%  "

# Row 21 - was Avalanche, now WrappedBTC
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "'26.662.00"
$ws.Range("E21").Value = "  -3.54%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.135"
$ws.Range("E22").Value = "  -3.42%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  -1.78%  "

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "'1.950.14"
$ws.Range("E24").Value = "  -7.62%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'151.76"
$ws.Range("E25").Value = "  -0.05%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "'1.846"
$ws.Range("E26").Value = "  -2.60%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -2.70%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "'2.067"
$ws.Range("E28").Value = "  -3.57%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "'5.089"
$ws.Range("E29").Value = "  -2.94%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "'115.08"
$ws.Range("E30").Value = "  -2.13%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.08858"
$ws.Range("E31").Value = "  -1.50%  "

# Row 32 - HuobiToken
$ws.Range("D32").Value = "'2.971"
$ws.Range("E32").Value = "  +0.87%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'4.409"
$ws.Range("E33").Value = "  -3.10%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  -4.17%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "'0.7135"
$ws.Range("E35").Value = "  -5.74%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  -2.61%  "

# Row 37 - Hedera
$ws.Range("E37").Value = "  -1.87%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.01920"
$ws.Range("E38").Value = "  -1.62%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "'2.398"
$ws.Range("E39").Value = "  -0.21%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "'2.918"
$ws.Range("E40").Value = "  -2.29%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "'7.139"
$ws.Range("E41").Value = "  -2.48%  "

# Row 42 - TheSandbox
$ws.Range("D42").Value = "'0.5140"
$ws.Range("E42").Value = "  -3.47%  "

# Row 43 - Algorand
$ws.Range("E43").Value = "  -2.16%  "

# Row 44 - Aptos
$ws.Range("E44").Value = "  -4.13%  "

# Row 45 - Decentraland
$ws.Range("E45").Value = "  -2.32%  "

# Row 46 - PaxDollar
$ws.Range("E46").Value = "  -0.20%  "

# Row 47 - was Quant, now EnergySwap
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.06"
$ws.Range("E47").Value = "  -4.65%  "

# Row 48 - was EnergySwap, now Quant
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'102.69"
$ws.Range("E48").Value = "  -2.19%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "'0.06264"
$ws.Range("E49").Value = "  -0.51%  "

# Row 50 - NEARProtocol
$ws.Range("D50").Value = "'1.611"
$ws.Range("E50").Value = "  -3.80%  "

# Row 51 - Aave
$ws.Range("E51").Value = "  -2.63%  "
